$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Duplicate the formatting of the last existing table row (row 55) into the
# new row 56 by copying the row and inserting it with a shift-down, which
# preserves the exact cell styles used throughout the table.
$ws.Range("B55:F55").Copy()
$ws.Range("B56:F56").Insert(-4121)

# Fill in the data for Post 46 (Dining Philosopher problem).
$ws.Range("B56").Value = 46
$ws.Range("C56").Value = "Dining Philosopher problem | Operating System - M03 P08"
$ws.Range("D56").Value = 44171
$ws.Range("E56").Value = "https://programmingport.hashnode.dev/dining-philosopher-problem-or-operating-system-m03-p08"
$ws.Range("F56").Value = "https://dev.to/rahulmishra05/dining-philosopher-problem-operating-system-m03-p08-fa5"

# Expand the structured table (Table2) so it covers the new row.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("B10:F56"))

# Match the author's final selection state.
$ws.Activate()
[void]$ws.Range("E56").Select()
